# feat: add 2022-Q3 data
#
# 1) Insert a new worksheet "2022-Q3" right after "总计" (i.e. right before
#    the existing "2021-Q3" sheet), populate it with the new quarter's fund
#    data (header "基金规模" instead of "基金金额").
# 2) Update the "总计" summary sheet: push the existing two rows down one
#    slot and insert the new "2022-Q3" totals row at the top.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet, positioned before "2021-Q3"
# ---------------------------------------------------------------------
$existingQ3 = $wb.Worksheets.Item("2021-Q3")
$newSheet = $wb.Worksheets.Add($existingQ3)
$newSheet.Name = "2022-Q3"

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$headerRange = $newSheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Row 2 - 008763
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'008763"
$newSheet.Range("C2").Value = "天弘越南市场股票（QDII）A"
$newSheet.Range("D2").Value = "'20.44"
$newSheet.Range("E2").Value = "'90.19"
$newSheet.Range("F2").Value = "'6.44"
$newSheet.Range("G2").Value = "'1.3163"
$newSheet.Range("H2").Value = 2

# Row 3 - 008764
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'008764"
$newSheet.Range("C3").Value = "天弘越南市场股票（QDII）C"
$newSheet.Range("D3").Value = "'15.02"
$newSheet.Range("E3").Value = "'90.19"
$newSheet.Range("F3").Value = "'6.44"
$newSheet.Range("G3").Value = "'0.9673"
$newSheet.Range("H3").Value = 2

$indexRange = $newSheet.Range("A2:A3")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: shift rows 2-3 down to 3-4, then write the
#    new "2022-Q3" row at row 2.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# old row 3 (2021-Q2) -> row 4
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = $totalSheet.Range("B3").Value2
$totalSheet.Range("C4").Value = $totalSheet.Range("C3").Value2
$totalSheet.Range("D4").Value = $totalSheet.Range("D3").Value2

# A4 needs the same bold/centered/bordered look as A2/A3
$a4 = $totalSheet.Range("A4")
$a4.Font.Bold = $true
$a4.HorizontalAlignment = -4108
$a4.VerticalAlignment = -4160
$a4.Borders.LineStyle = 1

# old row 2 (2021-Q3) -> row 3
$totalSheet.Range("B3").Value = $totalSheet.Range("B2").Value2
$totalSheet.Range("C3").Value = $totalSheet.Range("C2").Value2
$totalSheet.Range("D3").Value = $totalSheet.Range("D2").Value2

# new row 2 (2022-Q3)
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 2.28

# Restore the originally-active tab ("2021-Q2") so the only observable
# change is the data itself, not which sheet is selected.
$wb.Worksheets.Item("2021-Q2").Activate()

